# Apply the "poster color dis" update:
#  - Column G (load counts) values are shifted up by one row (G2 now holds
#    what used to be in G3, etc. through G15 holding what used to be in G16).
#  - The old G16 / H16 cells are removed (row 16 no longer has G/H values).
#  - A new row 17 is added with a STDEV formula in E17.
#  - The sheet's active selection moves to G20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "load" values for G2:G15 (previous G3:G16 values, shifted up one row)
$newLoadValues = @(40, 51, 58, 34, 57, 72, 47, 57, 42, 54, 40, 53, 65, 61)

for ($i = 0; $i -lt $newLoadValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $newLoadValues[$i]
}

# Remove the old trailing G16 / H16 values (row 16 shrinks back to A:E)
$ws.Range("G16:H16").Clear() | Out-Null

# Add the new row 17 with the STDEV formula in E17
$ws.Range("E17").Formula = "=STDEV(E2, E16)"

# Update the saved selection to match the author's final cursor position
$ws.Range("G20").Select()
